# Refresh crypto price/volume data (and reorder a few rows by rank)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.241.32"
$ws.Range("E2").Value = "  +1.39%  "
$ws.Range("D3").Value = "2.386.48"
$ws.Range("E3").Value = "  +3.88%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.64"
$ws.Range("E5").Value = "  +0.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.69"
$ws.Range("E6").Value = "  +2.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.510"
$ws.Range("E7").Value = "  +0.70%  "
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.503"
$ws.Range("E9").Value = "  +2.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.39"
$ws.Range("E10").Value = "  -0.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0792"
$ws.Range("E11").Value = "  +0.66%  "
$ws.Range("B12").Value = "Chainlink"
$ws.Range("C12").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "18.58"
$ws.Range("E12").Value = "  -2.34%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.122"
$ws.Range("E13").Value = "  +2.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.80"
$ws.Range("E14").Value = "  +0.61%  "
$ws.Range("D15").Value = "2.752.85"
$ws.Range("E15").Value = "  +3.33%  "
$ws.Range("D16").Value = "2.378.57"
$ws.Range("E16").Value = "  +4.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.811"
$ws.Range("E17").Value = "  +4.07%  "
$ws.Range("D18").Value = "43.244.87"
$ws.Range("E18").Value = "  +1.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.33"
$ws.Range("E19").Value = "  +1.83%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.37"
$ws.Range("E20").Value = "  +6.42%  "
$ws.Range("E21").Value = "  +0.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.39"
$ws.Range("E22").Value = "  +1.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.52"
$ws.Range("E23").Value = "  +0.75%  "
$ws.Range("E24").Value = "  -1.46%  "
$ws.Range("E25").Value = "  +1.53%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.90"
$ws.Range("E27").Value = "  +2.64%  "
$ws.Range("E28").Value = "  -0.18%  "
$ws.Range("E29").Value = "  +1.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.72"
$ws.Range("E30").Value = "  -0.27%  "
$ws.Range("E31").Value = "  +3.33%  "
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0729"
$ws.Range("E33").Value = "  +4.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.30"
$ws.Range("E34").Value = "  -1.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.87"
$ws.Range("E35").Value = "  +7.52%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.40"
$ws.Range("E36").Value = "  -0.42%  "
$ws.Range("E37").Value = "  -0.89%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.102"
$ws.Range("E38").Value = "  +1.93%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.83"
$ws.Range("E39").Value = "  +5.38%  "
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.80"
$ws.Range("E40").Value = "  +14.24%  "
$ws.Range("E41").Value = "  +0.71%  "
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "104.24"
$ws.Range("E42").Value = "  -37.05%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "1.948.81"
$ws.Range("E43").Value = "  -0.69%  "
$ws.Range("E44").Value = "  +1.11%  "
$ws.Range("E45").Value = "  +2.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.38"
$ws.Range("E46").Value = "  -10.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.77"
$ws.Range("E47").Value = "  +0.66%  "
$ws.Range("D48").Value = "2.605.20"
$ws.Range("E48").Value = "  +3.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.16"
$ws.Range("E49").Value = "  +0.32%  "
$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "72.42"
$ws.Range("E50").Value = "  +1.72%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.51"
$ws.Range("E51").Value = "  +1.96%  "
